# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from 2023-09-12 (45181) to 2023-09-13 (45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45182
}
